# Each study has been renamed to match the ones in the grant proposal, by
# introducing a new "study" column (the grant-proposal short name) right
# after the existing author/year column, which itself is relabeled
# "author_year". The plot was also made more compact by giving the new
# column a slightly narrower custom width than the old outcome column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column B; everything that used to be in B..K shifts
# one column right, into C..L.
$ws.Columns("B").Insert()

# Fill in the new "study" column (B) with the grant-proposal study names.
$ws.Range("B1").Value = "study"
# The newly inserted header cell should not carry the bold header style
# that the rest of row 1 has - reset it back to the default/normal style.
$ws.Range("B1").Style = "Normal"

$ws.Range("B2").Value = "CTSN Severe MR"
$ws.Range("B3").Value = "CTSN Moderate MR"
$ws.Range("B4").Value = "CTSN TR Trial"
$ws.Range("B5").Value = "CTSN AF Trial"

# The old "study" column (now column A) is relabeled "author_year".
$ws.Range("A1").Value = "author_year"

# Give the new study column a custom width (compacting the layout).
$ws.Columns("B").ColumnWidth = 23.42578125
